# Apply the "gh-pages output generated at 456a3b4" update to 杭州-漫展信息.xlsx
#
# Summary of the change:
#   * Sheet "展览": a new event ("杭州·HD·01") is inserted as row 42, pushing
#     the previous row 42 ("杭州·D3动漫游戏嘉年华") down to row 43 and the
#     previous row 43 ("杭州·理想乡动漫展-同人创作者大会") down to row 44.
#     The "想去人数" (F) counters also ticked up for many existing rows.
#   * Sheet "演出": two "想去人数" (F) counters ticked up.
#   * Sheet "全部类型": the same "想去人数" (F) counters ticked up across the
#     consolidated listing (no row insertion here - it already listed every
#     event, so only values change; note row 43 is not a "D3" row in this
#     sheet's numbering, it is the earlier "动漫作品主题音乐会" show, so all
#     edits are driven by explicit per-cell targets, never by row offsets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Insert a new blank row at position 42; this shifts the two rows that used
# to be at 42/43 down to 43/44 and carries their formatting with them.
$ws1.Rows.Item(42).Insert()

# The newly inserted row 42 only gets column A's formatting copied over (the
# insert left it using a slightly different auto-generated style); pull the
# canonical style from the row immediately below, which still carries the
# original "index column" formatting.
$ws1.Range("A43").Copy()
$ws1.Range("A42").PasteSpecial(-4122)

# New row 42: 杭州·HD·01 (the new event pushed onto the list). Column A is
# left exactly as the insert left it (value 41, same as the source data's
# "index" column never renumbers on insert in the original export).
$ws1.Range("A42").Value = 41

# B42 looks like a plain ISO date ("2024-08-17"); Excel's COM layer will
# happily reinterpret that as a real date serial unless we tell it up front
# that this cell is text. Force text entry (apostrophe prefix), then reset
# the cell style to the sheet's default so no stray number-format sticks
# around on it.
$ws1.Range("B42").Value = "'2024-08-17"
$ws1.Range("B42").Style = "Normal"

$ws1.Range("C42").Value = "杭州·HD·01"
$ws1.Range("D42").Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws1.Range("E42").Value = "2024.08.17 09:30-08.18 17:00"
$ws1.Range("F42").Value = 9
$ws1.Range("G42").Value = 75
$ws1.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=86332"
$ws1.Range("I42").Value = "//i2.hdslb.com/bfs/openplatform/202405/GBMur4hT1716145118862.jpeg"

# Row 43 (previously row 42: 杭州·D3动漫游戏嘉年华) - only the "想去人数"
# ticked up. The Insert() shifted the old row 42's index value (41) down
# into A43 along with everything else, but the source data keeps A43's
# original value (42) unchanged, so restore it explicitly.
$ws1.Range("A43").Value = 42
$ws1.Range("F43").Value = 310

# Row 44 (previously row 43: 杭州·理想乡动漫展-同人创作者大会) - new ordinal id
# and "想去人数" ticked up
$ws1.Range("A44").Value = 43
$ws1.Range("F44").Value = 3555

# Across-the-board "想去人数" (F column) refresh for the rest of sheet "展览"
$sheet1Counts = @{
    "F2"  = 573
    "F3"  = 5405
    "F7"  = 1006
    "F8"  = 375
    "F9"  = 1342
    "F12" = 3063
    "F13" = 1905
    "F14" = 118
    "F16" = 188
    "F17" = 17
    "F18" = 135
    "F20" = 969
    "F21" = 347
    "F22" = 47
    "F23" = 3518
    "F24" = 1106
    "F25" = 2792
    "F26" = 279
    "F27" = 1959
    "F28" = 4027
    "F29" = 107
    "F30" = 913
    "F31" = 461
    "F32" = 1279
    "F33" = 11
    "F35" = 994
    "F36" = 1262
    "F37" = 58
    "F38" = 1023
    "F39" = 662
    "F40" = 516
    "F41" = 403
}
foreach ($addr in $sheet1Counts.Keys) {
    $ws1.Range($addr).Value = $sheet1Counts[$addr]
}

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$sheet2Counts = @{
    "F22" = 36
    "F25" = 6
}
foreach ($addr in $sheet2Counts.Keys) {
    $ws2.Range($addr).Value = $sheet2Counts[$addr]
}

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Counts = @{
    "F2"  = 573
    "F3"  = 573
    "F4"  = 5407
    "F9"  = 375
    "F10" = 1342
    "F11" = 3063
    "F13" = 1905
    "F14" = 118
    "F17" = 188
    "F21" = 135
    "F22" = 969
    "F23" = 347
    "F24" = 3519
    "F27" = 1106
    "F28" = 2792
    "F29" = 1959
    "F30" = 4027
    "F32" = 107
    "F33" = 913
    "F34" = 1279
    "F36" = 994
    "F38" = 1262
    "F39" = 58
    "F40" = 1023
    "F42" = 662
    "F44" = 403
    "F45" = 36
    "F47" = 6
    "F48" = 310
    "F49" = 3555
}
foreach ($addr in $sheet4Counts.Keys) {
    $ws4.Range($addr).Value = $sheet4Counts[$addr]
}
